$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id test"
$ws.Range("A2").Value = "test"

$ws.Range("A5").Select()
